$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: bold a literal marker (e.g. "(a)") inside a paragraph's Range,
# searching forward from a given character offset within that paragraph's
# text. Returns the offset (within the paragraph text) just past the marker,
# so callers can chain successive searches left-to-right.
# ---------------------------------------------------------------------------
function Bold-Marker($paraIndex, $marker, $searchFrom) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $t = $rng.Text
    $pos = $t.IndexOf($marker, $searchFrom)
    if ($pos -lt 0) { throw "marker '$marker' not found in paragraph $paraIndex" }
    $s = $rng.Start + $pos
    $e = $s + $marker.Length
    $sub = $d.Range($s, $e)
    $sub.Bold = 1
    return $pos + $marker.Length
}

# 1) Question 2 (c): bold the "(c)" label.
Bold-Marker 19 "(c)" 0 | Out-Null

# 2) Question 3 (a)/(b)/(c): bold each label in turn.
$next = Bold-Marker 22 "(a)" 0
$next = Bold-Marker 22 "(b)" $next
$next = Bold-Marker 22 "(c)" $next

# 3) Relocate the (hidden) "_GoBack" bookmark so it wraps the bold "(c)"
#    label plus the single space that follows it, right before "cite the
#    research..." in that same paragraph.
$p22 = $d.Paragraphs.Item(22)
$rng22 = $p22.Range
$t22 = $rng22.Text
$cPos = $t22.IndexOf("(c)", 150)
$bkStart = $rng22.Start + $cPos
$bkEnd = $bkStart + 4
$bkRange = $d.Range($bkStart, $bkEnd)
$d.Bookmarks.Add("_GoBack", $bkRange)

# ---------------------------------------------------------------------------
# 4) Replace the trailing (now bookmark-less) empty paragraph at the very end
#    of the body with a blank paragraph, a page break, a bold
#    "Literature Cited (for test days)" heading, and two trailing blank
#    paragraphs.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $lastPara.Range.Duplicate
$insertPoint.Collapse(1)  # wdCollapseStart

$bodyFrag = '<w:p/>' + `
            '<w:p><w:r><w:br w:type="page"/></w:r></w:p>' + `
            '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Literature Cited (for test days)</w:t></w:r></w:p>' + `
            '<w:p/><w:p/>'
$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyFrag + '</w:body></w:document></pkg:xmlData>' + `
           '</pkg:part></pkg:package>'
$insertPoint.InsertXML($xmlFrag)

# ---------------------------------------------------------------------------
# 5) Refresh the cached PAGE-field text shown in the header from the stale
#    "4" to the correct "1".
# ---------------------------------------------------------------------------
$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdr.Range.Find.Execute("4", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2) | Out-Null

Write-Output "edit complete"
